$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "PI and hourly wages" label text (62% -> 65%)
$ws.Range("A23").Value = "PI and hourly wages (65%)"

# Update the PI summer salary request percentage (I12): 50% -> 100%
$ws.Range("I12").Value = 1

# Update the base PI salary figures for each year (B3:E3)
$ws.Range("B3").Formula = "=24061*I12"
$ws.Range("C3").Formula = "=24782.83*I12"
$ws.Range("D3").Formula = "=25526.31*I12"
$ws.Range("E3").Formula = "=26292.1*I12"

# Move the active cell selection to G20
$ws.Range("G20").Select() | Out-Null

$wb.Application.Calculate() | Out-Null
